$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 49987.25
$ws.Range("I86").Value = 63316.332
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 63316.332
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -62193.332
$ws.Range("N86").Value = -12246

$ws.Range("H87").Value = 60000
$ws.Range("J87").Value = 60000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62496

$ws.Range("H89").Value = 49987.25
$ws.Range("I89").Value = 63316.332
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 316581.66
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -310965.66
$ws.Range("N89").Value = -61232

$ws.Range("H90").Value = 60000
$ws.Range("J90").Value = 60000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -192480

$ws.Range("H112").Value = 2172.818
$ws.Range("I112").Value = 1400.3334
$ws.Range("J112").Value = 2462.5
$ws.Range("K112").Value = 4201.0002
$ws.Range("L112").Value = 7387.5
$ws.Range("M112").Value = -3093.0002
$ws.Range("N112").Value = -9603.5

$ws.Range("H113").Value = 6874.125
$ws.Range("I113").Value = 4999.5
$ws.Range("J113").Value = 7499
$ws.Range("K113").Value = 4999.5
$ws.Range("L113").Value = 7499
$ws.Range("M113").Value = -1745.5
$ws.Range("N113").Value = -14007

$ws.Range("H132").Value = 1325.2858
$ws.Range("I132").Value = 1341.3636
$ws.Range("K132").Value = 4024.0908
$ws.Range("M132").Value = -1494.0908

$ws.Range("H141").Value = 2773.0312
$ws.Range("I141").Value = 1937.4445
$ws.Range("J141").Value = 3100
$ws.Range("K141").Value = 5812.333500000001
$ws.Range("L141").Value = 9300
$ws.Range("M141").Value = -632.3335000000006
$ws.Range("N141").Value = -19660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 732.3333
$ws.Range("I5").Value = 678.8
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 678.8
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -566.8
$ws.Range("N5").Value = -1224

$ws.Range("H76").Value = 32746
$ws.Range("J76").Value = 32746
$ws.Range("L76").Value = 32746
$ws.Range("N76").Value = -33422

$ws.Range("H79").Value = 32746
$ws.Range("J79").Value = 32746
$ws.Range("L79").Value = 32746
$ws.Range("N79").Value = -35086

$ws.Range("H88").Value = 1604.8667
$ws.Range("I88").Value = 1099
$ws.Range("J88").Value = 2047.5
$ws.Range("K88").Value = 1099
$ws.Range("L88").Value = 2047.5
$ws.Range("M88").Value = -693
$ws.Range("N88").Value = -2859.5

$ws.Range("H91").Value = 1604.8667
$ws.Range("I91").Value = 1099
$ws.Range("J91").Value = 2047.5
$ws.Range("K91").Value = 1099
$ws.Range("L91").Value = 2047.5
$ws.Range("M91").Value = 305
$ws.Range("N91").Value = -4855.5

$ws.Range("H132").Value = 2901.111
$ws.Range("I132").Value = 2763.75
$ws.Range("K132").Value = 8291.25
$ws.Range("M132").Value = -5761.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 732.3333
$ws.Range("I4").Value = 678.8
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 678.8
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -563.8
$ws.Range("N4").Value = -1230

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 178.63158
$ws.Range("I7").Value = 80
$ws.Range("J7").Value = 288.22223
$ws.Range("K7").Value = 80
$ws.Range("L7").Value = 288.22223
$ws.Range("M7").Value = 33
$ws.Range("N7").Value = -514.2222300000001

$ws.Range("H58").Value = 1591.8948
$ws.Range("I58").Value = 1694.75
$ws.Range("J58").Value = 1517.091
$ws.Range("K58").Value = 1694.75
$ws.Range("L58").Value = 1517.091
$ws.Range("M58").Value = -1491.75
$ws.Range("N58").Value = -1923.091

$ws.Range("H132").Value = 3055.5715
$ws.Range("I132").Value = 3055.5715
$ws.Range("K132").Value = 9166.7145
$ws.Range("M132").Value = -6636.7145

$ws.Range("H136").Value = 1591.8948
$ws.Range("I136").Value = 1694.75
$ws.Range("J136").Value = 1517.091
$ws.Range("K136").Value = 5084.25
$ws.Range("L136").Value = 4551.272999999999
$ws.Range("M136").Value = -2534.25
$ws.Range("N136").Value = -9651.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 7495
$ws.Range("J81").Value = 7495
$ws.Range("L81").Value = 22485
$ws.Range("N81").Value = -24731

$ws.Range("H84").Value = 7495
$ws.Range("J84").Value = 7495
$ws.Range("L84").Value = 67455
$ws.Range("N84").Value = -78687

$ws.Range("H86").Value = 1098.5
$ws.Range("I86").Value = 700
$ws.Range("J86").Value = 1497
$ws.Range("K86").Value = 2100
$ws.Range("L86").Value = 4491
$ws.Range("M86").Value = -914
$ws.Range("N86").Value = -6863

$ws.Range("H89").Value = 1098.5
$ws.Range("I89").Value = 700
$ws.Range("J89").Value = 1497
$ws.Range("K89").Value = 6300
$ws.Range("L89").Value = 13473
$ws.Range("M89").Value = -372
$ws.Range("N89").Value = -25329

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1336.3334
$ws.Range("J6").Value = 1754.5
$ws.Range("L6").Value = 1754.5
$ws.Range("N6").Value = -1980.5

$ws.Range("H7").Value = 500663
$ws.Range("I7").Value = 750000
$ws.Range("J7").Value = 334438.34
$ws.Range("K7").Value = 750000
$ws.Range("L7").Value = 334438.34
$ws.Range("M7").Value = -749888
$ws.Range("N7").Value = -334662.34

$ws.Range("H8").Value = 500663
$ws.Range("I8").Value = 750000
$ws.Range("J8").Value = 334438.34
$ws.Range("K8").Value = 750000
$ws.Range("L8").Value = 334438.34
$ws.Range("M8").Value = -749861
$ws.Range("N8").Value = -334716.34

$ws.Range("H16").Value = 1336.3334
$ws.Range("J16").Value = 1754.5
$ws.Range("L16").Value = 1754.5
$ws.Range("N16").Value = -2254.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5315.1665
$ws.Range("I7").Value = 5288.2
$ws.Range("J7").Value = 5450
$ws.Range("K7").Value = 5288.2
$ws.Range("L7").Value = 5450
$ws.Range("M7").Value = -5176.2
$ws.Range("N7").Value = -5674

$ws.Range("H14").Value = 9666.666999999999
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 9666.666999999999
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 9666.666999999999
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -10010.667

$ws.Range("H16").Value = 2600
$ws.Range("J16").Value = 1200
$ws.Range("L16").Value = 1200
$ws.Range("N16").Value = -1540

$ws.Range("H93").Value = 2332.6667
$ws.Range("I93").Value = 1999.2
$ws.Range("K93").Value = 1999.2
$ws.Range("M93").Value = -751.2

$ws.Range("H126").Value = 5315.1665
$ws.Range("I126").Value = 5288.2
$ws.Range("J126").Value = 5450
$ws.Range("K126").Value = 15864.6
$ws.Range("L126").Value = 16350
$ws.Range("M126").Value = -13394.6
$ws.Range("N126").Value = -21290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 30125.334
$ws.Range("I41").Value = 29999
$ws.Range("J41").Value = 30188.5
$ws.Range("K41").Value = 29999
$ws.Range("L41").Value = 30188.5
$ws.Range("N41").Value = -30968.5
$ws.Range("M41").Value = -29609

$ws.Range("H74").Value = 16545.25
$ws.Range("I74").Value = 14500
$ws.Range("K74").Value = 14500
$ws.Range("M74").Value = -13564

$ws.Range("H77").Value = 16545.25
$ws.Range("I77").Value = 14500
$ws.Range("K77").Value = 43500
$ws.Range("M77").Value = -38820

$ws.Range("H136").Value = 3445.64
$ws.Range("I136").Value = 3244.2354
$ws.Range("J136").Value = 3873.625
$ws.Range("K136").Value = 9732.706200000001
$ws.Range("L136").Value = 11620.875
$ws.Range("M136").Value = -7182.706200000001
$ws.Range("N136").Value = -16720.875
